# Extend the rolling "forecasts_table" on both sheets (cases, deaths) with two more
# observation-date columns (AC, AD = 2020-05-08 / 2020-05-09) and two more forecast-origin
# rows (41, 42 = 2020-05-22 / 2020-05-23), matching the new diagonal of forecast values.
$wb = $excel.ActiveWorkbook

# ---- Sheet "cases" ----
$ws1 = $wb.Worksheets.Item("cases")

# Extend header row with the next two observation dates (AC1, AD1)
$ws1.Cells.Item(1, 29).NumberFormat = "@"
$ws1.Cells.Item(1, 29).Value = "2020-05-08"
$ws1.Cells.Item(1, 29).Style = "Normal"
$ws1.Cells.Item(1, 30).NumberFormat = "@"
$ws1.Cells.Item(1, 30).Value = "2020-05-09"
$ws1.Cells.Item(1, 30).Style = "Normal"

# Rows 2-26: widen the table with empty AC/AD cells
for ($r = 2; $r -le 26; $r++) {
    $ws1.Cells.Item($r, 29).Style = "Normal"
    $ws1.Cells.Item($r, 30).Style = "Normal"
}

# Rows 27-40: fill in forecast values / widen with empty cells as needed
$ws1.Cells.Item(27, 2).Value = 15741
$ws1.Cells.Item(27, 29).Style = "Normal"
$ws1.Cells.Item(27, 30).Style = "Normal"
$ws1.Cells.Item(28, 2).Value = 16929
$ws1.Cells.Item(28, 29).Value = 16412
$ws1.Cells.Item(28, 30).Style = "Normal"
$ws1.Cells.Item(29, 29).Value = 17636
$ws1.Cells.Item(29, 30).Value = 18141
$ws1.Cells.Item(30, 29).Value = 18410
$ws1.Cells.Item(30, 30).Value = 19124
$ws1.Cells.Item(31, 29).Value = 19168
$ws1.Cells.Item(31, 30).Value = 20087
$ws1.Cells.Item(32, 29).Value = 19955
$ws1.Cells.Item(32, 30).Value = 21083
$ws1.Cells.Item(33, 29).Value = 20766
$ws1.Cells.Item(33, 30).Value = 22054
$ws1.Cells.Item(34, 29).Value = 21578
$ws1.Cells.Item(34, 30).Value = 22837
$ws1.Cells.Item(35, 29).Value = 22456
$ws1.Cells.Item(35, 30).Value = 23764
$ws1.Cells.Item(36, 29).Value = 23032
$ws1.Cells.Item(36, 30).Value = 24519
$ws1.Cells.Item(37, 29).Value = 23386
$ws1.Cells.Item(37, 30).Value = 25119
$ws1.Cells.Item(38, 29).Value = 23921
$ws1.Cells.Item(38, 30).Value = 25746
$ws1.Cells.Item(39, 29).Value = 24518
$ws1.Cells.Item(39, 30).Value = 26362
$ws1.Cells.Item(40, 29).Value = 25451
$ws1.Cells.Item(40, 30).Value = 27208

# Rows 41-42: brand-new forecast-origin rows (2020-05-22, 2020-05-23)
$ws1.Cells.Item(41, 1).NumberFormat = "@"
$ws1.Cells.Item(41, 1).Value = "2020-05-22"
$ws1.Cells.Item(41, 1).Style = "Normal"
for ($c = 2; $c -le 28; $c++) {
    $ws1.Cells.Item(41, $c).Style = "Normal"
}
$ws1.Cells.Item(41, 29).Value = 26057
$ws1.Cells.Item(41, 30).Value = 27875

$ws1.Cells.Item(42, 1).NumberFormat = "@"
$ws1.Cells.Item(42, 1).Value = "2020-05-23"
$ws1.Cells.Item(42, 1).Style = "Normal"
for ($c = 2; $c -le 28; $c++) {
    $ws1.Cells.Item(42, $c).Style = "Normal"
}
$ws1.Cells.Item(42, 29).Style = "Normal"
$ws1.Cells.Item(42, 30).Value = 28615

# ---- Sheet "deaths" ----
$ws2 = $wb.Worksheets.Item("deaths")

# Extend header row with the next two observation dates (AC1, AD1)
$ws2.Cells.Item(1, 29).NumberFormat = "@"
$ws2.Cells.Item(1, 29).Value = "2020-05-08"
$ws2.Cells.Item(1, 29).Style = "Normal"
$ws2.Cells.Item(1, 30).NumberFormat = "@"
$ws2.Cells.Item(1, 30).Value = "2020-05-09"
$ws2.Cells.Item(1, 30).Style = "Normal"

# Rows 2-26: widen the table with empty AC/AD cells
for ($r = 2; $r -le 26; $r++) {
    $ws2.Cells.Item($r, 29).Style = "Normal"
    $ws2.Cells.Item($r, 30).Style = "Normal"
}

# Rows 27-40: fill in forecast values / widen with empty cells as needed
$ws2.Cells.Item(27, 2).Value = 1503
$ws2.Cells.Item(27, 29).Style = "Normal"
$ws2.Cells.Item(27, 30).Style = "Normal"
$ws2.Cells.Item(28, 2).Value = 1653
$ws2.Cells.Item(28, 29).Value = 1607
$ws2.Cells.Item(28, 30).Style = "Normal"
$ws2.Cells.Item(29, 29).Value = 1752
$ws2.Cells.Item(29, 30).Value = 1795
$ws2.Cells.Item(30, 29).Value = 1841
$ws2.Cells.Item(30, 30).Value = 1905
$ws2.Cells.Item(31, 29).Value = 1947
$ws2.Cells.Item(31, 30).Value = 2015
$ws2.Cells.Item(32, 29).Value = 2054
$ws2.Cells.Item(32, 30).Value = 2129
$ws2.Cells.Item(33, 29).Value = 2164
$ws2.Cells.Item(33, 30).Value = 2249
$ws2.Cells.Item(34, 29).Value = 2251
$ws2.Cells.Item(34, 30).Value = 2371
$ws2.Cells.Item(35, 29).Value = 2360
$ws2.Cells.Item(35, 30).Value = 2500
$ws2.Cells.Item(36, 29).Value = 2427
$ws2.Cells.Item(36, 30).Value = 2599
$ws2.Cells.Item(37, 29).Value = 2482
$ws2.Cells.Item(37, 30).Value = 2670
$ws2.Cells.Item(38, 29).Value = 2560
$ws2.Cells.Item(38, 30).Value = 2752
$ws2.Cells.Item(39, 29).Value = 2651
$ws2.Cells.Item(39, 30).Value = 2844
$ws2.Cells.Item(40, 29).Value = 2778
$ws2.Cells.Item(40, 30).Value = 2973

# Rows 41-42: brand-new forecast-origin rows (2020-05-22, 2020-05-23)
$ws2.Cells.Item(41, 1).NumberFormat = "@"
$ws2.Cells.Item(41, 1).Value = "2020-05-22"
$ws2.Cells.Item(41, 1).Style = "Normal"
for ($c = 2; $c -le 28; $c++) {
    $ws2.Cells.Item(41, $c).Style = "Normal"
}
$ws2.Cells.Item(41, 29).Value = 2862
$ws2.Cells.Item(41, 30).Value = 3081

$ws2.Cells.Item(42, 1).NumberFormat = "@"
$ws2.Cells.Item(42, 1).Value = "2020-05-23"
$ws2.Cells.Item(42, 1).Style = "Normal"
for ($c = 2; $c -le 28; $c++) {
    $ws2.Cells.Item(42, $c).Style = "Normal"
}
$ws2.Cells.Item(42, 29).Style = "Normal"
$ws2.Cells.Item(42, 30).Value = 3193
